$wb = $excel.ActiveWorkbook

# --- DegreeRequirement sheet: rename headers, add a computed "opening" column ---
$ws = $wb.Worksheets.Item("DegreeRequirement")

# Set C1 before A1 so new shared-string entries are minted in the same order
# as the authored workbook (RequirementID, then DegreeRequirementID).
$ws.Range("C1").Value = "RequirementID"
$ws.Range("A1").Value = "DegreeRequirementID"
$ws.Range("B1").Value = "DegreeID"

# New column D: one CONCATENATE formula per data row (entered per-cell so
# Excel does not collapse them into a shared-formula group).
for ($r = 2; $r -le 13; $r++) {
    $ws.Range("D$r").Formula = '=CONCATENATE("new DegreeRequirement{","DegreeRequirementID=",A:A,",","DegreeID=",B:B,",","RequirementID=",C:C,"},")'
}

$ws.Columns.Item(4).ColumnWidth = 80.65

# Make this the active sheet/tab and select the new column's data.
$ws.Activate() | Out-Null
$ws.Range("D2:D13").Select() | Out-Null
